$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("E16").Value = 7

$ws.Range("E17").Value = 3
$ws.Range("F17").Value = "I had stuff here but people recklessly reverted commits and then commited those reverts to the head branch. They did so without telling anyone and now I cannot figure out what I did specifically or what is missing from the application!"

$ws.Range("E18").Value = 1
$ws.Range("F18").Value = "Added new ""Likes"" table wrote backend code to update likes in the database."

$ws.Range("F19").Select()
